# EPBDS-10212 @class appears in REST response
# Adds a new "HomeCat extends Cat" datatype block (with a "smart" typed
# field) to the Rules sheet, mirroring the existing Datatype blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New field row - Boolean smart
$ws.Range("C40").Value = "smart"
$ws.Range("B40").Value = "Boolean"
$ws.Range("B40:C40").VerticalAlignment = -4108
$ws.Range("B40:C40").WrapText = $true

# New header row - "Datatype HomeCat extends Cat" (merged B39:C39, same
# look & feel as the other "Datatype ..." header rows, e.g. B10:C10).
$ws.Range("B39:C39").Merge() | Out-Null
$ws.Range("B39:C39").HorizontalAlignment = -4108
$ws.Range("B39:C39").VerticalAlignment = -4108
$ws.Range("B39:C39").WrapText = $true
$ws.Range("B39").Value = "Datatype HomeCat extends Cat"

$ws.Range("B50").Select() | Out-Null
